$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header cells: _old -> _FV2310, _new -> _FV2404
$oldSuffixCols = @(1,2,3,4,5,6,7,8,9,10)   # A..J
$newSuffixCols = @(12,13,14,15,16,17,18,19,20,21) # L..U

$baseNames = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Cells.Item(1, $oldSuffixCols[$i]).Value = $baseNames[$i] + "_FV2310"
    $ws.Cells.Item(1, $newSuffixCols[$i]).Value = $baseNames[$i] + "_FV2404"
}

# Freeze the top row (pane split after row 1)
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Convert the range into an Excel Table (ListObject)
$tableRange = $ws.Range("A1:U61")
$listObject = $ws.ListObjects.Add(1, $tableRange, [System.Type]::Missing, 1)
$listObject.Name = "Table1"
$listObject.TableStyle = ""
